# "Working on the Alt versions"
#
# 1) Duplicate the "Original" sheet, rename the copy to "Alt" and drop it at
#    the end of the tab strip.
# 2) Strip the picture/shape objects that came along with the copy (the Alt
#    sheet does not carry the header image).
# 3) Add a MAP/LAMBDA/LET based decrypt formula in A14 that spills down
#    through A19 (the alternate single-formula approach being explored).
# 4) MySingleFunction stops being the active sheet (Alt takes over) - update
#    its selection, and turn the "Answer Expected" header link in I1 into a
#    real hyperlink (pointing at the LinkedIn post), which also changes its
#    cell style to the builtin Hyperlink style.

$wb = $excel.ActiveWorkbook

# --- 1) Copy "Original" -> "Alt" (placed after the last existing sheet) ---
$wsOriginal = $wb.Worksheets.Item("Original")
$wsLastExisting = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsOriginal.Copy($null, $wsLastExisting)

$wsAlt = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsAlt.Name = "Alt"

# --- 2) Remove the copied-over picture/shape objects ---
for ($i = $wsAlt.Shapes.Count; $i -ge 1; $i--) {
    $wsAlt.Shapes.Item($i).Delete()
}

# --- 3) New MAP/LAMBDA/LET decrypt formula, spilling A14:A19 ---
$wsAlt.Activate()
$wsAlt.Range("A14").Formula2 = "=MAP(A2:A7,B2:B7,LAMBDA(a,k,LET(s,SEQUENCE(LEN(a)),c,CODE(MID(a,s,1)),f,FLOOR(c,32),CONCAT(CHAR(MOD(c-MID(REPT(k,9),s,1)-f,26)+f)))))"
$wsAlt.Range("C10").Select()

# --- 4) MySingleFunction: no longer the active tab, move selection ---
$wsFunc = $wb.Worksheets.Item("MySingleFunction")
$wsFunc.Activate()
$wsFunc.Range("A42").Select()
$wsFunc.Hyperlinks.Add($wsFunc.Range("I1"), "https://www.linkedin.com/feed/update/urn:li:activity:7214477472195592193/") | Out-Null

# --- Make "Alt" the active tab, matching the saved workbook view ---
$wsAlt.Activate()

Write-Output "done"
